$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A52").Value = "Emanuele Tomasoni "
$ws.Range("B52").Value = "MATTEO PILATI | Pinguini Trentini"
$ws.Range("C52").Value = "Elia Tomasoni | Demobusters"
$ws.Range("D52").Value = "Marco  Sartorelli | Modium"
$ws.Range("E52").Value = "Nicolò Speziali | F.C. Gorillaz"
$ws.Range("F52").Value = "Niccolò Orsi | SBARX"
